# Propensity Score Matching bullet:
#   ": Responsible for building the PySpark and SQL pipelines for data
#    pre-processing, ..."
# becomes
#   ": Responsible for building ETL pipelines with PySpark and SQL for data
#    pre-processing, ..."
#
# The "PySpark" run (wrapped in proofErr spell-check tags) and the
# " and SQL" run must stay completely untouched, so every edit below is
# scoped to stay strictly inside its own run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1 - split the run ": Responsible for building the " into three
# runs: ": Responsible for building" / " ETL pipelines with" / " "
# ---------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute(": Responsible for building the ", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
$rng.Text = ": Responsible for building"
$b1 = $rng.End

$ins1 = $d.Range($b1, $b1)
$ins1.InsertAfter(" ETL pipelines with")
$b2 = $ins1.End

$ins2 = $d.Range($b2, $b2)
$ins2.InsertAfter(" ")
$b3 = $ins2.End

# Re-assert the two new run boundaries (a plain text write coalesces
# adjacent same-formatted runs back together, so give each freshly
# inserted piece a tiny, reverted formatting nudge - this is enough to
# make it keep its own run instead of being re-merged with its
# neighbour).
$p1 = $d.Range($b1, $b2)
$p1.Bold = 1
$p1.Bold = 0

$p2 = $d.Range($b2, $b3)
$p2.Bold = 1
$p2.Bold = 0

# ---------------------------------------------------------------------
# Part 2 - drop the word "pipelines" from the final run of the bullet
# ---------------------------------------------------------------------
$rng2 = $d.Content
$null = $rng2.Find.Execute(" pipelines for data pre-processing, feature engineering, feature selection, and propensity score matching to reduce bias in observational data.   ", `
                            $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s = $rng2.Start
$e = $rng2.End

# The run to the left (" and SQL") has identical run formatting, so it
# would normally be coalesced into the same run as our edit (and its
# run-level rsid would leak onto the new text). Give it a momentary
# formatting nudge so it is left alone by the coalesce, then edit the
# target run, then revert the nudge as the very last touch on that
# boundary so " and SQL" ends up exactly as it started.
$leftNeighbor = $d.Range($s - 8, $s)
$leftNeighbor.Bold = 1

$rng3 = $d.Range($s, $e)
$rng3.Text = " for data pre-processing, feature engineering, feature selection, and propensity score matching to reduce bias in observational data.   "

$leftNeighbor2 = $d.Range($s - 8, $s)
$leftNeighbor2.Bold = 0
